$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("anotacoes")

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "27/12/2025 00:51"
$ws.Range("C9").Value = 542
$ws.Range("D9").Value = "Conhecimentos Específicos"
$ws.Range("E9").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("F9").Value = "Estudar TPM"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "27/12/2025 00:52"
$ws.Range("C10").Value = 542
$ws.Range("D10").Value = "Conhecimentos Específicos"
$ws.Range("E10").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("F10").Value = "Estudar esse tal de 8 S"
